{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces each two-digit-by-two-digit multiplication equation in the\n// worksheet table with a newly generated equation, in document order.\n//\n// The body of this file is the body of `async (context) => { ... }`.\n\nconst replacements = [\n  [\"26\u00d720=520\", \"59\u00d780=4720\"],\n  [\"16\u00d749=784\", \"93\u00d743=3999\"],\n  [\"29\u00d781=2349\", \"96\u00d794=9024\"],\n  [\"32\u00d773=2336\", \"93\u00d735=3255\"],\n  [\"49\u00d722=1078\", \"34\u00d750=1700\"],\n  [\"19\u00d768=1292\", \"87\u00d756=4872\"],\n  [\"20\u00d791=1820\", \"73\u00d794=6862\"],\n  [\"31\u00d711=341\", \"37\u00d713=481\"],\n  [\"85\u00d725=2125\", \"60\u00d712=720\"],\n  [\"18\u00d742=756\", \"83\u00d723=1909\"],\n  [\"42\u00d717=714\", \"16\u00d723=368\"],\n  [\"40\u00d723=920\", \"11\u00d796=1056\"],\n  [\"65\u00d728=1820\", \"38\u00d750=1900\"],\n  [\"71\u00d743=3053\", \"83\u00d721=1743\"],\n  [\"40\u00d782=3280\", \"47\u00d772=3384\"],\n  [\"61\u00d798=5978\", \"74\u00d721=1554\"],\n  [\"66\u00d717=1122\", \"99\u00d776=7524\"],\n  [\"94\u00d789=8366\", \"58\u00d742=2436\"],\n  [\"60\u00d774=4440\", \"93\u00d775=6975\"],\n  [\"48\u00d782=3936\", \"30\u00d740=1200\"],\n  [\"49\u00d742=2058\", \"78\u00d793=7254\"],\n  [\"26\u00d795=2470\", \"90\u00d756=5040\"],\n  [\"87\u00d771=6177\", \"81\u00d724=1944\"],\n  [\"61\u00d764=3904\", \"58\u00d713=754\"],\n  [\"24\u00d745=1080\", \"88\u00d770=6160\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  // Only the first occurrence should exist (each equation is unique), but\n  // guard against duplicates by only touching the first match.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Replaces each two-digit-by-two-digit multiplication equation in the\n# worksheet table with a newly generated equation, in document order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"26\u00d720=520\";  New = \"59\u00d780=4720\" },\n    @{ Old = \"16\u00d749=784\";  New = \"93\u00d743=3999\" },\n    @{ Old = \"29\u00d781=2349\"; New = \"96\u00d794=9024\" },\n    @{ Old = \"32\u00d773=2336\"; New = \"93\u00d735=3255\" },\n    @{ Old = \"49\u00d722=1078\"; New = \"34\u00d750=1700\" },\n    @{ Old = \"19\u00d768=1292\"; New = \"87\u00d756=4872\" },\n    @{ Old = \"20\u00d791=1820\"; New = \"73\u00d794=6862\" },\n    @{ Old = \"31\u00d711=341\";  New = \"37\u00d713=481\" },\n    @{ Old = \"85\u00d725=2125\"; New = \"60\u00d712=720\" },\n    @{ Old = \"18\u00d742=756\";  New = \"83\u00d723=1909\" },\n    @{ Old = \"42\u00d717=714\";  New = \"16\u00d723=368\" },\n    @{ Old = \"40\u00d723=920\";  New = \"11\u00d796=1056\" },\n    @{ Old = \"65\u00d728=1820\"; New = \"38\u00d750=1900\" },\n    @{ Old = \"71\u00d743=3053\"; New = \"83\u00d721=1743\" },\n    @{ Old = \"40\u00d782=3280\"; New = \"47\u00d772=3384\" },\n    @{ Old = \"61\u00d798=5978\"; New = \"74\u00d721=1554\" },\n    @{ Old = \"66\u00d717=1122\"; New = \"99\u00d776=7524\" },\n    @{ Old = \"94\u00d789=8366\"; New = \"58\u00d742=2436\" },\n    @{ Old = \"60\u00d774=4440\"; New = \"93\u00d775=6975\" },\n    @{ Old = \"48\u00d782=3936\"; New = \"30\u00d740=1200\" },\n    @{ Old = \"49\u00d742=2058\"; New = \"78\u00d793=7254\" },\n    @{ Old = \"26\u00d795=2470\"; New = \"90\u00d756=5040\" },\n    @{ Old = \"87\u00d771=6177\"; New = \"81\u00d724=1944\" },\n    @{ Old = \"61\u00d764=3904\"; New = \"58\u00d713=754\" },\n    @{ Old = \"24\u00d745=1080\"; New = \"88\u00d770=6160\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: could not find '$($pair.Old)'\"\n    }\n}\n"}
